# Merge the split "<id>...</id>" runs (e.g. "<id>" + "p14" + "6" + "v_1" +
# "</id>") back into a single run per occurrence, as produced by a fresh
# tc/tcn/tl download. Word's Find & Replace on the concatenated visible
# text collapses every run it touches into one run carrying the
# formatting of the first run matched, which is exactly the merge we need.

$d = $word.ActiveDocument

$replacements = @(
    "<id>p146v_1</id>",
    "<id>p147v_1</id>",
    "<id>p147v_2</id>",
    "<id>p147v_3</id>"
)

foreach ($text in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2)
}
